{"js": "const body = context.document.body;\n\n// Whole-body range: replacing it swaps out the entire body content (the\n// lone \"aaaaa\" paragraph) for the corrected paragraphs below, same as the\n// COM script replacing $d.Content.\nconst targetRange = body.getRange();\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p/>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Os componentes do sistema </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>vao</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>linkeditar</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> as classes que tem o mesmo </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>esteri\u00f3tipo</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> de camada da arquitetura</w:t></w:r>\n            <w:r><w:t>.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(ooxml, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p/>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">Os componentes do sistema </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>vao</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>linkeditar</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> as classes que tem o mesmo </w:t></w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r><w:t>esteri\u00f3tipo</w:t></w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r><w:t xml:space=\"preserve\"> de camada da arquitetura</w:t></w:r>\n            <w:r><w:t>.</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$d.Content.InsertXML($xml)\n"}
